$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.076.26"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.906.01"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7429"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.97%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.81"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.05%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3086"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -3.11%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.43"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -5.84%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06970"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  +0.35%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7664"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.899.83"
$ws.Range("E13").Value = "  -1.77%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.311"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -1.78%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.13"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -1.02%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "30.088.31"
$ws.Range("E17").Value = "  -0.85%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.070"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -1.00%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007823"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -2.26%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.90"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -4.97%  "
$ws.Range("D21").Value = "2.214.35"
$ws.Range("E21").Value = "  +1.20%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = $style
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +0.02%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.141"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +6.47%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.377"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.11%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.00"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +1.04%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.99"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.60%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1271"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -2.61%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.046"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -7.07%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.351"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.542"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.24%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.331"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -2.53%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.079"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.70%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05227"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -1.36%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.301"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -2.70%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7467"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -1.52%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -2.45%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01968"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +0.28%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.10%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.334"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -2.71%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4487"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -0.77%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.31"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -6.04%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.970"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  +0.06%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8400"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.02%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.727"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -0.10%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.80"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +0.12%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.920"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "2.090.30"
$ws.Range("E49").Value = "  +0.07%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.68"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -2.62%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1180"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -5.35%  "
